# Update the "想去人数" (want-to-go count) values in column F across the
# workbook's sheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1292
$ws1.Range("F4").Value = 935
$ws1.Range("F6").Value = 1740
$ws1.Range("F7").Value = 385
$ws1.Range("F8").Value = 1161
$ws1.Range("F10").Value = 7
$ws1.Range("F11").Value = 117
$ws1.Range("F16").Value = 141
$ws1.Range("F25").Value = 141
$ws1.Range("F28").Value = 302
$ws1.Range("F29").Value = 126

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 614

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 300

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 300
$ws4.Range("F4").Value = 1292
$ws4.Range("F5").Value = 935
$ws4.Range("F7").Value = 1740
$ws4.Range("F8").Value = 385
$ws4.Range("F9").Value = 1161
$ws4.Range("F12").Value = 7
$ws4.Range("F13").Value = 117
$ws4.Range("F18").Value = 141
$ws4.Range("F33").Value = 141
$ws4.Range("F36").Value = 302
$ws4.Range("F39").Value = 126
$ws4.Range("F42").Value = 614
